$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 226, shifting existing rows 226:261 down to 227:262
$ws.Rows.Item(226).Insert()

# Populate the newly inserted row 226 with the new data record
$ws.Cells.Item(226, 1).Value = 3
$ws.Cells.Item(226, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(226, 3).Value = "Coquimbo"
$ws.Cells.Item(226, 4).Value = 44984
$ws.Cells.Item(226, 5).Value = 5
$ws.Cells.Item(226, 6).Value = 100112030
$ws.Cells.Item(226, 7).Value = "Poroto granado"
$ws.Cells.Item(226, 8).Value = "Sin especificar"
$ws.Cells.Item(226, 9).Value = "Primera"
$ws.Cells.Item(226, 10).Value = 73
$ws.Cells.Item(226, 11).Value = 39000
$ws.Cells.Item(226, 12).Value = 40000
$ws.Cells.Item(226, 13).Value = 39521
$ws.Cells.Item(226, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(226, 15).Value = "Región de Valparaíso"
$ws.Cells.Item(226, 16).Value = 1581
$ws.Cells.Item(226, 17).Value = 25
$ws.Cells.Item(226, 18).Value = "Hortaliza"

# Keep the date-formatted style consistent with the other date cells in column D
$ws.Cells.Item(226, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
